$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed cryptocurrency price / 1h-volume figures, and restored the
# PancakeSwap / BabyDogeCoin row order (rows 47-48) to match the latest feed.
$updates = @{
    "D2" = "29.734.80"
    "E2" = "  -2.51%  "
    "D3" = "2.097.21"
    "E3" = "  -1.86%  "
    "D4" = "1.010"
    "E4" = "  +0.21%  "
    "D5" = "344.11"
    "E5" = "  -2.31%  "
    "E6" = "  +0.19%  "
    "D7" = "0.5167"
    "E7" = "  -1.64%  "
    "D8" = "0.4393"
    "E8" = "  -3.48%  "
    "D9" = "52.91"
    "E9" = "  -1.39%  "
    "D10" = "0.09273"
    "E10" = "  +1.28%  "
    "D11" = "1.166"
    "E11" = "  -2.30%  "
    "D12" = "24.95"
    "D13" = "2.102.29"
    "E13" = "  -1.30%  "
    "D14" = "8.263"
    "E14" = "  +1.20%  "
    "D15" = "6.771"
    "E15" = "  -1.76%  "
    "D16" = "99.61"
    "E16" = "  -1.88%  "
    "D17" = "0.00001155"
    "E17" = "  -1.14%  "
    "E18" = "  +0.13%  "
    "D19" = "20.83"
    "E19" = "  +1.30%  "
    "D20" = "0.06649"
    "E20" = "  -0.99%  "
    "E21" = "  +0.19%  "
    "D22" = "6.204"
    "E22" = "  -2.72%  "
    "D23" = "29.778.35"
    "D24" = "12.50"
    "E24" = "  -2.81%  "
    "D25" = "2.319"
    "E25" = "  -2.42%  "
    "D26" = "2.351.12"
    "E26" = "  -1.38%  "
    "D27" = "21.96"
    "D28" = "2.521"
    "E28" = "  -3.21%  "
    "D29" = "161.40"
    "E29" = "  -2.21%  "
    "D30" = "133.14"
    "E30" = "  -2.03%  "
    "D31" = "1.143"
    "E31" = "  -6.46%  "
    "E32" = "  -3.04%  "
    "D33" = "1.653"
    "E33" = "  -3.94%  "
    "D34" = "6.179"
    "E34" = "  -3.43%  "
    "D35" = "3.939"
    "E35" = "  -2.32%  "
    "D36" = "6.295"
    "E36" = "  +2.29%  "
    "D37" = "10.25"
    "E37" = "  -1.86%  "
    "D38" = "0.02579"
    "E38" = "  -2.53%  "
    "D39" = "0.06728"
    "E39" = "  -3.67%  "
    "D40" = "12.47"
    "E40" = "  -2.18%  "
    "D41" = "0.6898"
    "E41" = "  -1.41%  "
    "D42" = "0.2235"
    "E42" = "  -5.22%  "
    "D43" = "1.314"
    "E43" = "  +3.12%  "
    "D44" = "0.6798"
    "E44" = "  +4.23%  "
    "D45" = "14.30"
    "E45" = "  -3.15%  "
    "D46" = "2.325"
    "E46" = "  -1.23%  "
    "B47" = "PancakeSwap"
    "C47" = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
    "D47" = "3.629"
    "E47" = "  -2.59%  "
    "B48" = "BabyDogeCoin"
    "C48" = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
    "D48" = "0.00000000357"
    "E48" = "  -5.34%  "
    "D49" = "1.220"
    "E49" = "  -2.24%  "
    "D50" = "82.19"
    "E50" = "  -1.98%  "
    "E51" = "  -2.07%  "
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    # Force text so numeric-looking strings (e.g. "1.010", "29.734.80")
    # keep their exact literal form instead of being parsed as numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
    # Drop the temporary text format so the cell style matches the original
    # (no explicit style index), leaving only the value changed.
    $cell.ClearFormats()
}
